# Add a new "Notifications" row to the "Functionality (Client App)" sheet,
# with a status of "DEVELOPING" styled with the built-in "Calculation" cell
# style (orange bold text on light-grey fill with a thin grey border),
# centered horizontally and vertically - matching the other Status cells'
# center alignment. Finally move the active selection to H11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 19: Function / Status columns (A and C), matching the pattern of
# the existing rows (A has wrap text, C holds the status badge style).
$ws.Range("A19").Value = "Notifications"
$ws.Range("A19").WrapText = $true

$ws.Range("C19").Value = "DEVELOPING"
$ws.Range("C19").Style = "Calculation"
$ws.Range("C19").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C19").VerticalAlignment = -4108    # xlCenter

# Match the recorded selection left behind in the workbook after editing.
$ws.Range("H11").Select() | Out-Null
